$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "323.94"
Set-TextValue "E2" "-2.03%"
Set-TextValue "D3" "39.24"
Set-TextValue "E3" "-2.97%"
Set-TextValue "D4" "5.687"
Set-TextValue "E4" "7.49%"
Set-TextValue "D5" "0.08007"
Set-TextValue "E5" "-1.08%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "2.001"
Set-TextValue "E6" "3.90%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "4.494"
Set-TextValue "E7" "-0.59%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "8.625"
Set-TextValue "E8" "0.18%"
Set-TextValue "D9" "2.944"
Set-TextValue "E9" "-1.11%"
Set-TextValue "D10" "0.9269"
Set-TextValue "E10" "-1.03%"
Set-TextValue "E11" "-7.33%"
Set-TextValue "D12" "0.1983"
Set-TextValue "E12" "0.87%"
Set-TextValue "D13" "8.690"
Set-TextValue "E13" "24.27%"
Set-TextValue "D14" "0.09184"
Set-TextValue "E14" "-0.61%"
Set-TextValue "E15" "0.93%"
Set-TextValue "D16" "0.1049"
Set-TextValue "E16" "9.56%"
Set-TextValue "D17" "0.001292"
Set-TextValue "E17" "-2.36%"
Set-TextValue "D18" "0.006130"
Set-TextValue "E18" "1.80%"
Set-TextValue "D20" "0.3473"
Set-TextValue "E20" "-1.37%"
Set-TextValue "D21" "0.1372"
Set-TextValue "E21" "3.67%"
Set-TextValue "D22" "0.2413"
Set-TextValue "E22" "-5.86%"
Set-TextValue "D23" "0.04408"
Set-TextValue "E23" "-0.23%"
Set-TextValue "D24" "0.001264"
Set-TextValue "E24" "3.58%"
Set-TextValue "D25" "0.004626"
Set-TextValue "E25" "7.46%"
Set-TextValue "E26" "-3.36%"
Set-TextValue "D39" "0.02493"
Set-TextValue "E39" "-0.26%"
Set-TextValue "D40" "0.05345"
Set-TextValue "E40" "3.61%"
Set-TextValue "D41" "0.007465"
Set-TextValue "E41" "-3.01%"
Set-TextValue "D42" "0.009616"
Set-TextValue "E42" "5.00%"
Set-TextValue "D43" "0.1405"
Set-TextValue "E43" "-1.54%"
Set-TextValue "E44" "-2.49%"
Set-TextValue "E45" "-2.08%"
Set-TextValue "D46" "0.00006719"
Set-TextValue "E46" "0.94%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "-0.05%"
Set-TextValue "D48" "0.002972"
Set-TextValue "E48" "-11.15%"
Set-TextValue "D49" "0.002292"
Set-TextValue "E49" "-4.62%"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "-0.05%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "-0.05%"
